$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra data rows (rows 3-9); everything below row 2 shifts up/out
$ws.Rows("3:9").Delete()

# Update header row (row 1) to standardized uppercase/underscore format
$ws.Range("A1").Value = "NETWORK"
$ws.Range("B1").Value = "DC_TYPE"
$ws.Range("C1").Value = "SUB_TYPE"
$ws.Range("D1").Value = "DC_NUMBER"
$ws.Range("E1").Value = "CITY"
$ws.Range("F1").Value = "STATE"

# Update the DC Number value on the remaining data row (kept as text, like the
# rest of this numeric-looking-but-text column) and strip the temporary
# formatting back off so no extra styling is left behind on the cell
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "6011"
$ws.Range("D2").ClearFormats()
